# Change "en_cn" to "en_us" in the shared string used by mini_poker!B2
# and move the active cell selection from D29 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mini_poker")

$ws.Range("B2").Value = "en_us"

$ws.Activate()
$ws.Range("B3").Select()
